$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 for the newest week's record; this shifts the
# existing rows 2:47 down to 3:48 (dates/quality/volume/prices all move
# down one row, oldest-row-47 data becomes row 48) and extends the used
# range to A1:T48.
$ws.Rows.Item(2).Insert()

# The Insert() call copies the formatting of the row above (the bold
# header row) into the freshly inserted row; strip that back off so the
# new data row matches the plain (unstyled) look of every other data row.
$ws.Range("A2:T2").ClearFormats()

# Populate the newly inserted row 2 with this week's Chirimoya record.
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "Macroferia Regional de Talca"
$ws.Range("C2").Value = "Maule"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D2").Value = "2021-10-27"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100107
$ws.Range("H2").Value = "Otros"
$ws.Range("I2").Value = 100107002
$ws.Range("J2").Value = "Chirimoya"
$ws.Range("K2").Value = "Cultivar IV Región"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 25000
$ws.Range("O2").Value = 25000
$ws.Range("P2").Value = 25000
$ws.Range("Q2").Value = '$/bandeja 10 kilos'
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 2500
$ws.Range("T2").Value = 10
